# Indian MF 1st Stab - add 9 new weekly columns (Jun_16 .. Sep_08) with
# their analyst-rating note cells, matching the author's update.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert 9 new columns in front of the existing weekly columns (old B:V
#    shifts right by 9 to become K:AE).
# ---------------------------------------------------------------------------
$ws.Range("B1:J1").EntireColumn.Insert()

# ---------------------------------------------------------------------------
# 2. New week headers for the 9 inserted columns (most-recent-first, same
#    convention as the rest of row 1).
# ---------------------------------------------------------------------------
$ws.Range("B1").Value = "Sep_08"
$ws.Range("C1").Value = "Aug_25"
$ws.Range("D1").Value = "Aug_04"
$ws.Range("E1").Value = "Jul_23"
$ws.Range("F1").Value = "Jul_17"
$ws.Range("G1").Value = "Jul_07"
$ws.Range("H1").Value = "Jun_30"
$ws.Range("I1").Value = "Jun_24"
$ws.Range("J1").Value = "Jun_16"

# ---------------------------------------------------------------------------
# 3. The new columns behave like every other week column: default content is
#    the "UN" placeholder for every analyst-firm row (2-33).
# ---------------------------------------------------------------------------
$ws.Range("B2:J33").Value = "UN"

# ---------------------------------------------------------------------------
# 4. Drop in the actual rating-action notes that were published during the
#    new weeks, reusing the same green ("raises"/"initiates") and
#    pink ("downgrades"/"lowers") highlight colors already used elsewhere in
#    the sheet (copied verbatim from existing analogous cells V7 and G20).
# ---------------------------------------------------------------------------
$greenColor = $ws.Range("V7").Interior.Color
$pinkColor  = $ws.Range("G20").Interior.Color

# Needham & Company LLC (row 2) - Raises Target
$ws.Range("C2").Value = "8/19/2019,Raises Target,Buy,`$92.00"
$ws.Range("C2").Interior.Color = $greenColor

# Piper Jaffray Companies (row 10) - Raises Target
$ws.Range("C10").Value = "8/15/2019,Raises Target,Overweight,`$83.00 -> `$86.00"
$ws.Range("C10").Interior.Color = $greenColor

# Piper Jaffray Companies (row 10) - Initiates (no highlight color)
$ws.Range("J10").Value = "6/11/2019,Initiates,Overweight -> Overweight,`$83.00"

# ValuEngine (row 11) - Downgrades
$ws.Range("F11").Value = "7/12/2019,Downgrades,Buy -> Hold,"
$ws.Range("F11").Interior.Color = $pinkColor

# Bank of America (row 20) - Downgrades
$ws.Range("F20").Value = "7/16/2019,Downgrades,Buy -> Neutral,`$76.00"
$ws.Range("F20").Interior.Color = $pinkColor

# Barclays (row 27) - Downgrades
$ws.Range("F27").Value = "7/15/2019,Downgrades,Overweight -> Equal Weight,`$86.00 -> `$78.00"
$ws.Range("F27").Interior.Color = $pinkColor
